$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 3
    3 = 2
    4 = 2
    5 = 2
    6 = 1
    7 = 1
    8 = 1
    9 = 0
    10 = 1
    11 = 3
    12 = 2
    13 = 0
    14 = 1
    15 = 2
    16 = 2
    17 = 1
    18 = 1
    19 = 0
    20 = 2
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 2
    26 = 2
    27 = 1
    28 = 2
    29 = 1
    30 = 1
    31 = 2
    32 = 0
    33 = 2
    34 = 3
    35 = 1
    36 = 0
    37 = 1
    38 = 2
    39 = 2
    40 = 1
    41 = 2
    42 = 1
    43 = 1
    44 = 3
    45 = 3
    46 = 0
    47 = 0
    48 = 3
    49 = 1
    50 = 1
    51 = 3
    52 = 0
    53 = 2
    54 = 1
    55 = 1
    56 = 1
    57 = 2
    58 = 3
    59 = 2
    60 = 2
    61 = 2
    62 = 1
    64 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
